$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "28.418.74"
$ws.Range("E2").Value = "  -0.10%  "
Set-TextValue "D3" "1.822.34"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "314.55"
$ws.Range("E5").Value = "  -0.75%  "
Set-TextValue "D6" "1.002"
$ws.Range("E6").Value = "  -0.02%  "
Set-TextValue "D7" "0.5117"
$ws.Range("E7").Value = "  -3.54%  "
Set-TextValue "D8" "0.3928"
$ws.Range("E8").Value = "  -3.40%  "
Set-TextValue "D9" "0.07660"
$ws.Range("E9").Value = "  +1.26%  "
Set-TextValue "D10" "41.69"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  +0.60%  "
Set-TextValue "D13" "6.267"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  -0.02%  "
Set-TextValue "D15" "7.486"
$ws.Range("E15").Value = "  -1.39%  "
Set-TextValue "D16" "1.823.68"
$ws.Range("E16").Value = "  -0.63%  "
Set-TextValue "D17" "93.07"
$ws.Range("E17").Value = "  +3.86%  "
Set-TextValue "D18" "0.00001098"
$ws.Range("E18").Value = "  +2.24%  "
Set-TextValue "D19" "0.06642"
$ws.Range("E19").Value = "  +0.46%  "
Set-TextValue "D20" "17.70"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  -0.09%  "
Set-TextValue "D22" "6.098"
$ws.Range("E22").Value = "  +0.48%  "
Set-TextValue "D23" "28.431.25"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  -1.44%  "
Set-TextValue "D25" "2.257"
$ws.Range("E25").Value = "  +6.52%  "
Set-TextValue "D26" "20.82"
$ws.Range("E26").Value = "  +1.15%  "
Set-TextValue "D27" "156.00"
$ws.Range("E27").Value = "  -0.76%  "
Set-TextValue "D28" "2.033.48"
$ws.Range("E28").Value = "  -0.74%  "
Set-TextValue "D29" "2.385"
$ws.Range("E29").Value = "  -3.74%  "
Set-TextValue "D30" "123.81"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D31" "0.1097"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.107"
$ws.Range("E32").Value = "  -1.66%  "
Set-TextValue "D33" "5.646"
$ws.Range("E33").Value = "  -0.75%  "
Set-TextValue "D34" "3.659"
$ws.Range("E34").Value = "  +0.02%  "
Set-TextValue "D35" "0.07072"
$ws.Range("E35").Value = "  -1.61%  "
Set-TextValue "D36" "0.2208"
Set-TextValue "D37" "0.02327"
$ws.Range("E37").Value = "  -0.77%  "
Set-TextValue "D38" "5.164"
$ws.Range("E38").Value = "  -1.97%  "
Set-TextValue "D39" "8.751"
$ws.Range("E39").Value = "  -0.54%  "
Set-TextValue "D40" "0.6251"
$ws.Range("E41").Value = "  -1.45%  "
Set-TextValue "D42" "1.173"
$ws.Range("E42").Value = "  -1.72%  "
Set-TextValue "D43" "1.000"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  -1.41%  "
Set-TextValue "D45" "13.35"
$ws.Range("E45").Value = "  -1.06%  "
Set-TextValue "D46" "3.727"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +0.33%  "
Set-TextValue "D48" "124.33"
$ws.Range("E48").Value = "  -1.44%  "
Set-TextValue "D49" "1.977"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -0.08%  "
Set-TextValue "D51" "0.06895"
$ws.Range("E51").Value = "  -0.12%  "